# Auto-generated edit script: refresh market-price derived columns (H:N)
# for the Leve profit tables across sheets, per scheduled-runner update.
$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:L126").ClearContents()
$ws.Range("H127:M127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("H129:M129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:L133").ClearContents()
$ws.Range("H134:L134").ClearContents()
$ws.Range("H135:M135").ClearContents()
$ws.Range("H136:L136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:L139").ClearContents()
$ws.Range("H140:L140").ClearContents()
$ws.Range("H141:M141").ClearContents()

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 25000
$ws.Range("J101").Value = 25000
$ws.Range("L101").Value = 25000
$ws.Range("N101").Value = -31490

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 24800
$ws.Range("I97").Value = 24800
$ws.Range("K97").Value = 24800
$ws.Range("M97").Value = -23809
$ws.Range("H103").Value = 15828.5
$ws.Range("J103").Value = 15828.5
$ws.Range("L103").Value = 15828.5
$ws.Range("N103").Value = -18172.5
$ws.Range("H134").Value = 7882.2
$ws.Range("I134").Value = 7882.2
$ws.Range("K134").Value = 23646.6
$ws.Range("M134").Value = -21111.6

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1652.8334
$ws.Range("I132").Value = 980.5
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 2941.5
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -411.5
$ws.Range("N132").Value = -14052.5

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3747.875
$ws.Range("J80").Value = 3747.875
$ws.Range("L80").Value = 11243.625
$ws.Range("N80").Value = -13115.625
$ws.Range("H83").Value = 3747.875
$ws.Range("J83").Value = 3747.875
$ws.Range("L83").Value = 33730.875
$ws.Range("N83").Value = -43090.875
$ws.Range("H113").Value = 690.6
$ws.Range("I113").Value = 76.5
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 229.5
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = 1940.5
$ws.Range("N113").Value = -7640
$ws.Range("N137").ClearContents()
$ws.Range("H137").Value = 1500
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125:L125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("H129:L129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:L133").ClearContents()
$ws.Range("H134:L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H136:L136").ClearContents()
$ws.Range("H137:L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138:L138").ClearContents()
$ws.Range("H139:L139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H140:L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

